$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old F column ("ppm") contained duplicate/incorrect data. D1 header was
# blank even though D already held the proper "ppm" values. Fix the header
# and remove the stray F column, shifting sample_size/t_results/significance
# one column to the left.
$ws.Range("D1").Value = "ppm"
$ws.Columns("F").Delete()
